$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove rows 85-91 (data shrank from 91 to 84 rows)
$ws.Rows("85:91").Delete() | Out-Null

# Update data values for rows 3-84 (A,B,C columns)
$ws.Cells.Item(3, 1).Value = [double]"0"
$ws.Cells.Item(3, 2).Value = [double]"0.003164556962025316"
$ws.Cells.Item(3, 3).Value = [double]"0.9999995459485628"
$ws.Cells.Item(4, 1).Value = [double]"0"
$ws.Cells.Item(4, 2).Value = [double]"0.1075949367088608"
$ws.Cells.Item(4, 3).Value = [double]"0.9999837135908592"
$ws.Cells.Item(5, 1).Value = [double]"0"
$ws.Cells.Item(5, 2).Value = [double]"0.1645569620253164"
$ws.Cells.Item(5, 3).Value = [double]"0.9999819518009659"
$ws.Cells.Item(6, 1).Value = [double]"0"
$ws.Cells.Item(6, 2).Value = [double]"0.1993670886075949"
$ws.Cells.Item(6, 3).Value = [double]"0.9999758764822407"
$ws.Cells.Item(7, 1).Value = [double]"0"
$ws.Cells.Item(7, 2).Value = [double]"0.2626582278481013"
$ws.Cells.Item(7, 3).Value = [double]"0.999973429149117"
$ws.Cells.Item(8, 1).Value = [double]"0"
$ws.Cells.Item(8, 2).Value = [double]"0.3417721518987342"
$ws.Cells.Item(8, 3).Value = [double]"0.9999574050658746"
$ws.Cells.Item(9, 1).Value = [double]"0"
$ws.Cells.Item(9, 2).Value = [double]"0.4715189873417722"
$ws.Cells.Item(9, 3).Value = [double]"0.9999568472928365"
$ws.Cells.Item(10, 1).Value = [double]"0"
$ws.Cells.Item(10, 2).Value = [double]"0.5348101265822784"
$ws.Cells.Item(10, 3).Value = [double]"0.9999329170756178"
$ws.Cells.Item(11, 1).Value = [double]"0"
$ws.Cells.Item(11, 2).Value = [double]"0.5411392405063291"
$ws.Cells.Item(11, 3).Value = [double]"0.9999307273736462"
$ws.Cells.Item(12, 1).Value = [double]"0"
$ws.Cells.Item(12, 2).Value = [double]"0.5981012658227848"
$ws.Cells.Item(12, 3).Value = [double]"0.9999033681477413"
$ws.Cells.Item(13, 1).Value = [double]"0"
$ws.Cells.Item(13, 2).Value = [double]"0.6044303797468354"
$ws.Cells.Item(13, 3).Value = [double]"0.999898361402715"
$ws.Cells.Item(14, 1).Value = [double]"0"
$ws.Cells.Item(14, 2).Value = [double]"0.819620253164557"
$ws.Cells.Item(14, 3).Value = [double]"0.9988632513207271"
$ws.Cells.Item(15, 1).Value = [double]"0"
$ws.Cells.Item(15, 2).Value = [double]"0.8259493670886076"
$ws.Cells.Item(15, 3).Value = [double]"0.9988259859885452"
$ws.Cells.Item(16, 1).Value = [double]"0"
$ws.Cells.Item(16, 2).Value = [double]"0.9525316455696202"
$ws.Cells.Item(16, 3).Value = [double]"0.92043335169492"
$ws.Cells.Item(17, 1).Value = [double]"0.002590673575129534"
$ws.Cells.Item(17, 2).Value = [double]"0.9525316455696202"
$ws.Cells.Item(17, 3).Value = [double]"0.9155338700793891"
$ws.Cells.Item(18, 1).Value = [double]"0.002590673575129534"
$ws.Cells.Item(18, 2).Value = [double]"0.9715189873417721"
$ws.Cells.Item(18, 3).Value = [double]"0.7970923412247629"
$ws.Cells.Item(19, 1).Value = [double]"0.005181347150259068"
$ws.Cells.Item(19, 2).Value = [double]"0.9715189873417721"
$ws.Cells.Item(19, 3).Value = [double]"0.7572346215442045"
$ws.Cells.Item(20, 1).Value = [double]"0.005181347150259068"
$ws.Cells.Item(20, 2).Value = [double]"0.9746835443037974"
$ws.Cells.Item(20, 3).Value = [double]"0.7333479396452087"
$ws.Cells.Item(21, 1).Value = [double]"0.007772020725388601"
$ws.Cells.Item(21, 2).Value = [double]"0.9746835443037974"
$ws.Cells.Item(21, 3).Value = [double]"0.5643853474357633"
$ws.Cells.Item(22, 1).Value = [double]"0.007772020725388601"
$ws.Cells.Item(22, 2).Value = [double]"0.9778481012658228"
$ws.Cells.Item(22, 3).Value = [double]"0.5615030152257449"
$ws.Cells.Item(23, 1).Value = [double]"0.01036269430051814"
$ws.Cells.Item(23, 2).Value = [double]"0.9778481012658228"
$ws.Cells.Item(23, 3).Value = [double]"0.5358458416298769"
$ws.Cells.Item(24, 1).Value = [double]"0.01036269430051814"
$ws.Cells.Item(24, 2).Value = [double]"0.9810126582278481"
$ws.Cells.Item(24, 3).Value = [double]"0.5016934764603386"
$ws.Cells.Item(25, 1).Value = [double]"0.0155440414507772"
$ws.Cells.Item(25, 2).Value = [double]"0.9810126582278481"
$ws.Cells.Item(25, 3).Value = [double]"0.4408519898001989"
$ws.Cells.Item(26, 1).Value = [double]"0.0155440414507772"
$ws.Cells.Item(26, 2).Value = [double]"0.9841772151898734"
$ws.Cells.Item(26, 3).Value = [double]"0.4257251893438878"
$ws.Cells.Item(27, 1).Value = [double]"0.01813471502590673"
$ws.Cells.Item(27, 2).Value = [double]"0.9841772151898734"
$ws.Cells.Item(27, 3).Value = [double]"0.4247332665424868"
$ws.Cells.Item(28, 1).Value = [double]"0.01813471502590673"
$ws.Cells.Item(28, 2).Value = [double]"0.9873417721518988"
$ws.Cells.Item(28, 3).Value = [double]"0.3993195514066594"
$ws.Cells.Item(29, 1).Value = [double]"0.0310880829015544"
$ws.Cells.Item(29, 2).Value = [double]"0.9873417721518988"
$ws.Cells.Item(29, 3).Value = [double]"0.1809462156984304"
$ws.Cells.Item(30, 1).Value = [double]"0.0310880829015544"
$ws.Cells.Item(30, 2).Value = [double]"0.990506329113924"
$ws.Cells.Item(30, 3).Value = [double]"0.1638741932635325"
$ws.Cells.Item(31, 1).Value = [double]"0.03626943005181347"
$ws.Cells.Item(31, 2).Value = [double]"0.990506329113924"
$ws.Cells.Item(31, 3).Value = [double]"0.1505566303638075"
$ws.Cells.Item(32, 1).Value = [double]"0.03626943005181347"
$ws.Cells.Item(32, 2).Value = [double]"0.9936708860759493"
$ws.Cells.Item(32, 3).Value = [double]"0.1311976253808664"
$ws.Cells.Item(33, 1).Value = [double]"0.05440414507772021"
$ws.Cells.Item(33, 2).Value = [double]"0.9936708860759493"
$ws.Cells.Item(33, 3).Value = [double]"0.09605094088768698"
$ws.Cells.Item(34, 1).Value = [double]"0.05440414507772021"
$ws.Cells.Item(34, 2).Value = [double]"0.9968354430379747"
$ws.Cells.Item(34, 3).Value = [double]"0.08879582777208314"
$ws.Cells.Item(35, 1).Value = [double]"0.1398963730569948"
$ws.Cells.Item(35, 2).Value = [double]"0.9968354430379747"
$ws.Cells.Item(35, 3).Value = [double]"0.005297595959677674"
$ws.Cells.Item(36, 1).Value = [double]"0.1398963730569948"
$ws.Cells.Item(36, 2).Value = [double]"1"
$ws.Cells.Item(36, 3).Value = [double]"0.00499298470230999"
$ws.Cells.Item(37, 1).Value = [double]"0.2590673575129533"
$ws.Cells.Item(37, 2).Value = [double]"1"
$ws.Cells.Item(37, 3).Value = [double]"0.0007587549035045091"
$ws.Cells.Item(38, 1).Value = [double]"0.2642487046632124"
$ws.Cells.Item(38, 2).Value = [double]"1"
$ws.Cells.Item(38, 3).Value = [double]"0.0007446243916517802"
$ws.Cells.Item(39, 1).Value = [double]"0.461139896373057"
$ws.Cells.Item(39, 2).Value = [double]"1"
$ws.Cells.Item(39, 3).Value = [double]"0.0001193154005394725"
$ws.Cells.Item(40, 1).Value = [double]"0.4663212435233161"
$ws.Cells.Item(40, 2).Value = [double]"1"
$ws.Cells.Item(40, 3).Value = [double]"0.000118972224571397"
$ws.Cells.Item(41, 1).Value = [double]"0.4740932642487047"
$ws.Cells.Item(41, 2).Value = [double]"1"
$ws.Cells.Item(41, 3).Value = [double]"0.0001166810563556537"
$ws.Cells.Item(42, 1).Value = [double]"0.4792746113989637"
$ws.Cells.Item(42, 2).Value = [double]"1"
$ws.Cells.Item(42, 3).Value = [double]"0.0001165202059854132"
$ws.Cells.Item(43, 1).Value = [double]"0.4870466321243523"
$ws.Cells.Item(43, 2).Value = [double]"1"
$ws.Cells.Item(43, 3).Value = [double]"0.0001157611042139578"
$ws.Cells.Item(44, 1).Value = [double]"0.4948186528497409"
$ws.Cells.Item(44, 2).Value = [double]"1"
$ws.Cells.Item(44, 3).Value = [double]"0.0001061165967368477"
$ws.Cells.Item(45, 1).Value = [double]"0.5"
$ws.Cells.Item(45, 2).Value = [double]"1"
$ws.Cells.Item(45, 3).Value = [double]"0.0001037658387108537"
$ws.Cells.Item(46, 1).Value = [double]"0.5854922279792746"
$ws.Cells.Item(46, 2).Value = [double]"1"
$ws.Cells.Item(46, 3).Value = [double]"6.865013877328436E-05"
$ws.Cells.Item(47, 1).Value = [double]"0.5906735751295337"
$ws.Cells.Item(47, 2).Value = [double]"1"
$ws.Cells.Item(47, 3).Value = [double]"6.741729851961142E-05"
$ws.Cells.Item(48, 1).Value = [double]"0.655440414507772"
$ws.Cells.Item(48, 2).Value = [double]"1"
$ws.Cells.Item(48, 3).Value = [double]"3.020863898602339E-05"
$ws.Cells.Item(49, 1).Value = [double]"0.6632124352331606"
$ws.Cells.Item(49, 2).Value = [double]"1"
$ws.Cells.Item(49, 3).Value = [double]"2.717523226163769E-05"
$ws.Cells.Item(50, 1).Value = [double]"0.6683937823834197"
$ws.Cells.Item(50, 2).Value = [double]"1"
$ws.Cells.Item(50, 3).Value = [double]"2.679913585586272E-05"
$ws.Cells.Item(51, 1).Value = [double]"0.6761658031088082"
$ws.Cells.Item(51, 2).Value = [double]"1"
$ws.Cells.Item(51, 3).Value = [double]"2.547910662230077E-05"
$ws.Cells.Item(52, 1).Value = [double]"0.6839378238341969"
$ws.Cells.Item(52, 2).Value = [double]"1"
$ws.Cells.Item(52, 3).Value = [double]"2.360123143643051E-05"
$ws.Cells.Item(53, 1).Value = [double]"0.689119170984456"
$ws.Cells.Item(53, 2).Value = [double]"1"
$ws.Cells.Item(53, 3).Value = [double]"2.233359883915264E-05"
$ws.Cells.Item(54, 1).Value = [double]"0.7046632124352331"
$ws.Cells.Item(54, 2).Value = [double]"1"
$ws.Cells.Item(54, 3).Value = [double]"2.064634759226348E-05"
$ws.Cells.Item(55, 1).Value = [double]"0.7098445595854922"
$ws.Cells.Item(55, 2).Value = [double]"1"
$ws.Cells.Item(55, 3).Value = [double]"2.022176663324144E-05"
$ws.Cells.Item(56, 1).Value = [double]"0.7150259067357513"
$ws.Cells.Item(56, 2).Value = [double]"1"
$ws.Cells.Item(56, 3).Value = [double]"1.958158128202159E-05"
$ws.Cells.Item(57, 1).Value = [double]"0.7383419689119171"
$ws.Cells.Item(57, 2).Value = [double]"1"
$ws.Cells.Item(57, 3).Value = [double]"1.942990376491338E-05"
$ws.Cells.Item(58, 1).Value = [double]"0.7435233160621761"
$ws.Cells.Item(58, 2).Value = [double]"1"
$ws.Cells.Item(58, 3).Value = [double]"1.938676330173595E-05"
$ws.Cells.Item(59, 1).Value = [double]"0.7538860103626943"
$ws.Cells.Item(59, 2).Value = [double]"1"
$ws.Cells.Item(59, 3).Value = [double]"1.904161405902729E-05"
$ws.Cells.Item(60, 1).Value = [double]"0.7564766839378239"
$ws.Cells.Item(60, 2).Value = [double]"1"
$ws.Cells.Item(60, 3).Value = [double]"1.852984187485917E-05"
$ws.Cells.Item(61, 1).Value = [double]"0.7616580310880829"
$ws.Cells.Item(61, 2).Value = [double]"1"
$ws.Cells.Item(61, 3).Value = [double]"1.850323660706329E-05"
$ws.Cells.Item(62, 1).Value = [double]"0.8160621761658031"
$ws.Cells.Item(62, 2).Value = [double]"1"
$ws.Cells.Item(62, 3).Value = [double]"1.838631090559549E-05"
$ws.Cells.Item(63, 1).Value = [double]"0.8264248704663213"
$ws.Cells.Item(63, 2).Value = [double]"1"
$ws.Cells.Item(63, 3).Value = [double]"1.834548750244648E-05"
$ws.Cells.Item(64, 1).Value = [double]"0.8316062176165803"
$ws.Cells.Item(64, 2).Value = [double]"1"
$ws.Cells.Item(64, 3).Value = [double]"1.800052091589528E-05"
$ws.Cells.Item(65, 1).Value = [double]"0.8367875647668394"
$ws.Cells.Item(65, 2).Value = [double]"1"
$ws.Cells.Item(65, 3).Value = [double]"1.798980141497424E-05"
$ws.Cells.Item(66, 1).Value = [double]"0.8393782383419689"
$ws.Cells.Item(66, 2).Value = [double]"1"
$ws.Cells.Item(66, 3).Value = [double]"1.796397152686283E-05"
$ws.Cells.Item(67, 1).Value = [double]"0.8497409326424871"
$ws.Cells.Item(67, 2).Value = [double]"1"
$ws.Cells.Item(67, 3).Value = [double]"1.7866460986933E-05"
$ws.Cells.Item(68, 1).Value = [double]"0.8523316062176166"
$ws.Cells.Item(68, 2).Value = [double]"1"
$ws.Cells.Item(68, 3).Value = [double]"1.78267917938149E-05"
$ws.Cells.Item(69, 1).Value = [double]"0.8601036269430051"
$ws.Cells.Item(69, 2).Value = [double]"1"
$ws.Cells.Item(69, 3).Value = [double]"1.759263585719613E-05"
$ws.Cells.Item(70, 1).Value = [double]"0.8626943005181347"
$ws.Cells.Item(70, 2).Value = [double]"1"
$ws.Cells.Item(70, 3).Value = [double]"1.743065504805421E-05"
$ws.Cells.Item(71, 1).Value = [double]"0.8782383419689119"
$ws.Cells.Item(71, 2).Value = [double]"1"
$ws.Cells.Item(71, 3).Value = [double]"1.69768737826589E-05"
$ws.Cells.Item(72, 1).Value = [double]"0.883419689119171"
$ws.Cells.Item(72, 2).Value = [double]"1"
$ws.Cells.Item(72, 3).Value = [double]"1.664772221791218E-05"
$ws.Cells.Item(73, 1).Value = [double]"0.8937823834196891"
$ws.Cells.Item(73, 2).Value = [double]"1"
$ws.Cells.Item(73, 3).Value = [double]"1.649444137575953E-05"
$ws.Cells.Item(74, 1).Value = [double]"0.9067357512953368"
$ws.Cells.Item(74, 2).Value = [double]"1"
$ws.Cells.Item(74, 3).Value = [double]"1.482443977866163E-05"
$ws.Cells.Item(75, 1).Value = [double]"0.9222797927461139"
$ws.Cells.Item(75, 2).Value = [double]"1"
$ws.Cells.Item(75, 3).Value = [double]"1.446977253421284E-05"
$ws.Cells.Item(76, 1).Value = [double]"0.9300518134715026"
$ws.Cells.Item(76, 2).Value = [double]"1"
$ws.Cells.Item(76, 3).Value = [double]"1.440529585541281E-05"
$ws.Cells.Item(77, 1).Value = [double]"0.9326424870466321"
$ws.Cells.Item(77, 2).Value = [double]"1"
$ws.Cells.Item(77, 3).Value = [double]"1.38887389686793E-05"
$ws.Cells.Item(78, 1).Value = [double]"0.9378238341968912"
$ws.Cells.Item(78, 2).Value = [double]"1"
$ws.Cells.Item(78, 3).Value = [double]"1.352421696406817E-05"
$ws.Cells.Item(79, 1).Value = [double]"0.9455958549222798"
$ws.Cells.Item(79, 2).Value = [double]"1"
$ws.Cells.Item(79, 3).Value = [double]"1.337292879253144E-05"
$ws.Cells.Item(80, 1).Value = [double]"0.9507772020725389"
$ws.Cells.Item(80, 2).Value = [double]"1"
$ws.Cells.Item(80, 3).Value = [double]"1.310152241684699E-05"
$ws.Cells.Item(81, 1).Value = [double]"0.9585492227979274"
$ws.Cells.Item(81, 2).Value = [double]"1"
$ws.Cells.Item(81, 3).Value = [double]"1.257812753688347E-05"
$ws.Cells.Item(82, 1).Value = [double]"0.9689119170984456"
$ws.Cells.Item(82, 2).Value = [double]"1"
$ws.Cells.Item(82, 3).Value = [double]"1.239076024698817E-05"
$ws.Cells.Item(83, 1).Value = [double]"0.9818652849740933"
$ws.Cells.Item(83, 2).Value = [double]"1"
$ws.Cells.Item(83, 3).Value = [double]"1.237930150123696E-05"
$ws.Cells.Item(84, 1).Value = [double]"1"
$ws.Cells.Item(84, 2).Value = [double]"1"
$ws.Cells.Item(84, 3).Value = [double]"1.620559177063308E-06"
